$d = $word.ActiveDocument

# --- Step 1: Update the paragraph describing model usage / evaluation metrics ---
# (Split " from the workshop, I tried every one ... " into several runs, adding the
#  new sentences about training/evaluating each model with the chosen metrics.)
$p10 = $d.Paragraphs.Item(10)

$para10Xml = @"
<?xml version="1.0" encoding="UTF-8" standalone="yes"?>
<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">
<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">
<pkg:xmlData>
<w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">
<w:body>
<w:p><w:pPr><w:spacing w:after="0"/></w:pPr><w:r><w:t>Using the seven models provided</w:t></w:r><w:r><w:t xml:space="preserve"> from the workshop, </w:t></w:r><w:r><w:t>I trained each model and evaluated them using Balanced Accuracy, F1 Macro</w:t></w:r><w:r><w:t xml:space="preserve"> and Prediction Time as the main metrics. Additionally, confusion matrices were analysed to identify False Positives</w:t></w:r><w:r><w:t>, False Alarms, and Gray Areas</w:t></w:r><w:r><w:t>.</w:t></w:r><w:r><w:t xml:space="preserve"> I also did some very light edits to some of the models </w:t></w:r><w:r><w:t>to help them perform better without using up too much of my time.</w:t></w:r></w:p>
</w:body>
</w:document>
</pkg:xmlData>
</pkg:part>
</pkg:package>
"@

$null = $p10.Range.InsertXML($para10Xml)

# --- Step 2: Remove the embedded model-performance graph image ---
if ($d.Shapes.Count -gt 0) {
    $d.Shapes.Item(1).Delete()
}

# --- Step 3: Replace the intro sentence of the exclusion paragraph ("Using a
#     balanced model performance graph..." -> "From the above criteria, ") and
#     append the new paragraphs describing the final model choice, the
#     ThresholdClassifier tuning process, and the overall conclusion. ---
$p12 = $d.Paragraphs.Item(12)
$paraCountBefore = $d.Paragraphs.Count
$wasFinalParagraph = ($p12.Range.End -eq $d.Content.End)

$para12PlusXml = @"
<?xml version="1.0" encoding="UTF-8" standalone="yes"?>
<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">
<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">
<pkg:xmlData>
<w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">
<w:body>
<w:p><w:pPr><w:spacing w:after="0"/></w:pPr><w:r><w:t>From the above criteria,</w:t></w:r><w:r><w:t xml:space="preserve"> </w:t></w:r><w:r><w:t xml:space="preserve">Neural Network (MLP) and k-Nearest Neighbours were the first 2 models I excluded. Their False Positive rates were too high to be considered. </w:t></w:r><w:r><w:t xml:space="preserve">SVM </w:t></w:r><w:r><w:t>still had comparatively high False Positive rates while Random Forest models were much slower (0.06s) compared to the other models (&lt;0.01s)</w:t></w:r><w:r><w:t xml:space="preserve">. I decided to keep </w:t></w:r><w:r><w:t>Decision Tree (for having the lowest False Positive rate) and Gradient Boosting</w:t></w:r><w:r><w:t xml:space="preserve"> (for its extremely low False Alarm rates) to do more detailed testing.</w:t></w:r></w:p><w:p><w:pPr><w:spacing w:after="0"/></w:pPr></w:p><w:p><w:pPr><w:spacing w:after="0"/></w:pPr><w:r><w:t>After additional testing, we chose Gradient Boosting as the final model for this datathon</w:t></w:r><w:r><w:t xml:space="preserve">. </w:t></w:r></w:p><w:p><w:pPr><w:spacing w:after="0"/></w:pPr></w:p><w:p><w:pPr><w:spacing w:after="0"/></w:pPr><w:r><w:t>This model is</w:t></w:r><w:r><w:t xml:space="preserve"> then</w:t></w:r><w:r><w:t xml:space="preserve"> constructed in a custom ThresholdClassifier to adjust </w:t></w:r><w:r><w:t xml:space="preserve">how the %confidence for a </w:t></w:r><w:r><w:t>sample</w:t></w:r><w:r><w:t xml:space="preserve"> to be placed at a certain </w:t></w:r><w:r><w:t>class</w:t></w:r><w:r><w:t>.</w:t></w:r><w:r><w:t xml:space="preserve"> This classifier also included a reject margin that reassigns uncertain samples to class 2 instead of risking a wrong class 1 or class 3 classification.</w:t></w:r><w:r><w:t xml:space="preserve"> By adding the reject margin, the model actually did slightly worse producing more false alarms</w:t></w:r><w:r><w:t xml:space="preserve"> (pred = 2, true = 1)</w:t></w:r><w:r><w:t xml:space="preserve"> but I felt this was important to</w:t></w:r><w:r><w:t xml:space="preserve"> keep false positives low.</w:t></w:r></w:p><w:p><w:pPr><w:spacing w:after="0"/></w:pPr></w:p><w:p><w:pPr><w:spacing w:after="0"/></w:pPr><w:r><w:t>The final Gradient Boosting model together with the ThresholdClassifier pipeline provided the best</w:t></w:r><w:r><w:t xml:space="preserve"> trade-off between having </w:t></w:r><w:r><w:t xml:space="preserve">very low false positive rate yet low enough false alarms </w:t></w:r><w:r><w:t xml:space="preserve">to ensure </w:t></w:r><w:r><w:t xml:space="preserve">trust in clinicians. This model is therefore suitable in helping obstetric practitioners in identifying </w:t></w:r><w:r><w:t>fetal cases at risk without generating excessive false alarms.</w:t></w:r></w:p>
</w:body>
</w:document>
</pkg:xmlData>
</pkg:part>
</pkg:package>
"@

$null = $p12.Range.InsertXML($para12PlusXml)

# InsertXML on the last paragraph of the body leaves one extra empty trailing
# paragraph behind (it keeps the original closing paragraph mark). Merge it away
# so the new last paragraph becomes the document's true final paragraph, matching
# the target structure (no stray empty paragraph before the section properties).
if ($wasFinalParagraph) {
    $countAfter = $d.Paragraphs.Count
    if ($countAfter -gt $paraCountBefore) {
        $secondLast = $d.Paragraphs.Item($countAfter - 1)
        $markRange = $d.Range($secondLast.Range.End - 1, $secondLast.Range.End)
        $markRange.Delete()
    }
}
